$wb = $excel.ActiveWorkbook

# zh-cn sheet: row for 8afa32e3-... file (row 3) - Correspond Handoff Datetime (E)
# and Correspond Handback DateTime (H) get refreshed to newer timestamps.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-23 10:48:38"
$wsZhCn.Range("H3").Value = "2016-03-23 10:49:03"

# de-de sheet: same row/columns updated with its own timestamps.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-23 10:48:42"
$wsDeDe.Range("H3").Value = "2016-03-23 10:49:10"
